# Reorders the data rows 15-38 on the active sheet according to the
# permutation observed in the target diff. Row positions (and therefore
# any row-level formatting) stay put; only the cell *contents* of each
# row move to a different row.
#
# target row -> source row (i.e. target row ends up holding what used
# to be in source row before this edit)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$mapping = @{
    15 = 25
    16 = 28
    17 = 29
    18 = 32
    19 = 23
    20 = 17
    21 = 19
    22 = 20
    23 = 31
    24 = 21
    25 = 15
    26 = 30
    27 = 34
    28 = 38
    29 = 33
    30 = 27
    31 = 18
    32 = 26
    33 = 24
    34 = 22
    35 = 37
    36 = 36
    37 = 16
    38 = 35
}

$firstRow = 15
$lastRow  = 38
$firstCol = "A"
$lastCol  = "AY"

# 1. Snapshot every source row's full A:AY contents BEFORE any writes,
#    since several rows are both a source and a destination.
$snapshots = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rangeAddr = "$firstCol$r`:$lastCol$r"
    $snapshots[$r] = $ws.Range($rangeAddr).Value2
}

# 2. Write each snapshot into its destination row.
for ($t = $firstRow; $t -le $lastRow; $t++) {
    $src = $mapping[$t]
    $destAddr = "$firstCol$t`:$lastCol$t"
    $ws.Range($destAddr).Value2 = $snapshots[$src]
}
